$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Selplg"
$ws.Range("C2").Value = "Selp"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 78.14530633333334
$ws.Range("H2").Value = 234.435919
$ws.Range("I2").Value = 0.9738103308619316
$ws.Range("J2").Value = 0.9738103308619316
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 17.47090666666666
$ws.Range("N2").Value = 52.41271999999999
$ws.Range("O2").Value = 0.9803569739482672
$ws.Range("P2").Value = 0.9803569739482673
$ws.Range("Q2").Value = 1365.269353387742
$ws.Range("R2").Value = 12287.42418048968
$ws.Range("S2").Value = 0.9546817491633641
$ws.Range("T2").Value = 0.9546817491633642

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Selplg"
$ws.Range("C3").Value = "Selp"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 78.14530633333334
$ws.Range("H3").Value = 234.435919
$ws.Range("I3").Value = 0.9738103308619316
$ws.Range("J3").Value = 0.9738103308619316
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.16459
$ws.Range("N3").Value = 0.49377
$ws.Range("O3").Value = 0.0092357516081294
$ws.Range("P3").Value = 0.0092357516081294
$ws.Range("Q3").Value = 12.86193596940333
$ws.Range("R3").Value = 115.75742372463
$ws.Range("S3").Value = 0.008993870329271108
$ws.Range("T3").Value = 0.008993870329271108

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Selplg"
$ws.Range("C4").Value = "Selp"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 78.14530633333334
$ws.Range("H4").Value = 234.435919
$ws.Range("I4").Value = 0.9738103308619316
$ws.Range("J4").Value = 0.9738103308619316
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1854676666666667
$ws.Range("N4").Value = 0.556403
$ws.Range("O4").Value = 0.01040727444360334
$ws.Range("P4").Value = 0.01040727444360334
$ws.Range("Q4").Value = 14.49342762659522
$ws.Range("R4").Value = 130.440848639357
$ws.Range("S4").Value = 0.0101347113692963
$ws.Range("T4").Value = 0.0101347113692963

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Selplg"
$ws.Range("C5").Value = "Selp"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9293790000000001
$ws.Range("H5").Value = 2.788137
$ws.Range("I5").Value = 0.01158148728249443
$ws.Range("J5").Value = 0.01158148728249443
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 17.47090666666666
$ws.Range("N5").Value = 52.41271999999999
$ws.Range("O5").Value = 0.9803569739482672
$ws.Range("P5").Value = 0.9803569739482673
$ws.Range("Q5").Value = 16.23709376696
$ws.Range("R5").Value = 146.13384390264
$ws.Range("S5").Value = 0.01135399182608658
$ws.Range("T5").Value = 0.01135399182608658

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Selplg"
$ws.Range("C6").Value = "Selp"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9293790000000001
$ws.Range("H6").Value = 2.788137
$ws.Range("I6").Value = 0.01158148728249443
$ws.Range("J6").Value = 0.01158148728249443
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.16459
$ws.Range("N6").Value = 0.49377
$ws.Range("O6").Value = 0.0092357516081294
$ws.Range("P6").Value = 0.0092357516081294
$ws.Range("Q6").Value = 0.15296648961
$ws.Range("R6").Value = 1.37669840649
$ws.Range("S6").Value = 0.0001069637397938281
$ws.Range("T6").Value = 0.0001069637397938281

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Selplg"
$ws.Range("C7").Value = "Selp"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9293790000000001
$ws.Range("H7").Value = 2.788137
$ws.Range("I7").Value = 0.01158148728249443
$ws.Range("J7").Value = 0.01158148728249443
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.1854676666666667
$ws.Range("N7").Value = 0.556403
$ws.Range("O7").Value = 0.01040727444360334
$ws.Range("P7").Value = 0.01040727444360334
$ws.Range("Q7").Value = 0.172369754579
$ws.Range("R7").Value = 1.551327791211
$ws.Range("S7").Value = 0.0001205317166140214
$ws.Range("T7").Value = 0.0001205317166140214

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Selplg"
$ws.Range("C8").Value = "Selp"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.172262
$ws.Range("H8").Value = 3.516786
$ws.Range("I8").Value = 0.01460818185557397
$ws.Range("J8").Value = 0.01460818185557397
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 17.47090666666666
$ws.Range("N8").Value = 52.41271999999999
$ws.Range("O8").Value = 0.9803569739482672
$ws.Range("P8").Value = 0.9803569739482673
$ws.Range("Q8").Value = 20.48047999088
$ws.Range("R8").Value = 184.32431991792
$ws.Range("S8").Value = 0.01432123295881648
$ws.Range("T8").Value = 0.01432123295881648

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Selplg"
$ws.Range("C9").Value = "Selp"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.172262
$ws.Range("H9").Value = 3.516786
$ws.Range("I9").Value = 0.01460818185557397
$ws.Range("J9").Value = 0.01460818185557397
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.16459
$ws.Range("N9").Value = 0.49377
$ws.Range("O9").Value = 0.0092357516081294
$ws.Range("P9").Value = 0.0092357516081294
$ws.Range("Q9").Value = 0.19294260258
$ws.Range("R9").Value = 1.73648342322
$ws.Range("S9").Value = 0.000134917539064464
$ws.Range("T9").Value = 0.000134917539064464

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Selplg"
$ws.Range("C10").Value = "Selp"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.172262
$ws.Range("H10").Value = 3.516786
$ws.Range("I10").Value = 0.01460818185557397
$ws.Range("J10").Value = 0.01460818185557397
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.1854676666666667
$ws.Range("N10").Value = 0.556403
$ws.Range("O10").Value = 0.01040727444360334
$ws.Range("P10").Value = 0.01040727444360334
$ws.Range("Q10").Value = 0.217416697862
$ws.Range("R10").Value = 1.956750280758
$ws.Range("S10").Value = 0.000152031357693025
$ws.Range("T10").Value = 0.000152031357693025
